# Regenerate orders with updated distance/size codes.
# The Condition / Filename_Left / Filename_Right / Distance / Size columns
# encode experiment parameters as substrings inside the cell text
# (e.g. "Face07_D80_S20", "Face07_D80_S20_l.png", "D80", "S30", ...).
# This commit renumbers the distance levels (D64->D69, D80->D86, D51->D55)
# and the "S30" size level to "S31" everywhere those tokens occur.
#
# Because the codes are embedded as substrings throughout many different
# strings, the simplest faithful way to reproduce the change is a global
# find & replace across the whole used range, same as a human editor would
# do in the Excel UI (Ctrl+H) before regenerating/saving the order file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# Order matters only in that none of the replacement targets collide with
# any of the other search tokens, so a straightforward sequential pass is
# sufficient and safe.
$used.Replace("D80", "D86") | Out-Null
$used.Replace("D64", "D69") | Out-Null
$used.Replace("D51", "D55") | Out-Null
$used.Replace("S30", "S31") | Out-Null
